# Resample the per-trial stimulus rows (rows 2-41) from the existing pool of
# 20 unique kitchens-block trials: keep the positional columns (subject_id,
# task, block_total, block_scene, trial_block, target_cat) as-is, bump the
# running trial_total counter (col F) to continue past the previous block,
# and re-draw the per-stimulus columns (category, cond_cat, correct_answer,
# stimulus, conceptual, perceptual, typicality, n, p_*, r_*) from another row
# in the same pool, per the mapping below.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow  = 41

# destination row -> source row (within the original, pre-edit data) that
# supplies the "drawn stimulus" columns for that destination row.
$map = @{
    2 = 20
    3 = 38
    4 = 32
    5 = 41
    6 = 17
    7 = 26
    8 = 31
    9 = 28
    10 = 24
    11 = 19
    12 = 39
    13 = 14
    14 = 36
    15 = 3
    16 = 5
    17 = 34
    18 = 11
    19 = 37
    20 = 30
    21 = 6
    22 = 7
    23 = 29
    24 = 35
    25 = 25
    26 = 21
    27 = 27
    28 = 16
    29 = 33
    30 = 18
    31 = 2
    32 = 40
    33 = 10
    34 = 13
    35 = 23
    36 = 8
    37 = 4
    38 = 12
    39 = 22
    40 = 15
    41 = 9
}

# Columns carried over from the source row: category, cond_cat,
# correct_answer, stimulus, conceptual, perceptual, typicality, n,
# p_typicality, p_conceptual, p_perceptual, r_typicality, r_conceptual,
# r_perceptual
$drawnCols = 8,9,11,12,13,14,15,16,17,18,19,20,21,22

# Snapshot the pre-edit values for the drawn columns of every row before
# writing anything back (since several rows draw from each other).
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowVals = @()
    foreach ($c in $drawnCols) {
        $rowVals += $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowVals
}

for ($dst = $firstRow; $dst -le $lastRow; $dst++) {
    # trial_total (col F) keeps counting up across blocks.
    $ws.Cells.Item($dst, 6).Value = $dst + 161

    $src = $map[$dst]
    $vals = $snapshot[$src]
    $i = 0
    foreach ($c in $drawnCols) {
        $ws.Cells.Item($dst, $c).Value = $vals[$i]
        $i++
    }
}
